$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45205
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 11000
$ws.Range("L2").Value = 12000
$ws.Range("M2").Value = 11500
$ws.Range("P2").Value = 639

# Row 3
$ws.Range("D3").Value = 45175
$ws.Range("J3").Value = 250
$ws.Range("K3").Value = 11000
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = 11500
$ws.Range("P3").Value = 639

# Row 5
$ws.Range("D5").Value = 45092
$ws.Range("J5").Value = 210
$ws.Range("K5").Value = 10000
$ws.Range("L5").Value = 11000
$ws.Range("M5").Value = 10714
$ws.Range("P5").Value = 595

# Row 7
$ws.Range("D7").Value = 44792
$ws.Range("J7").Value = 160
$ws.Range("K7").Value = 9000
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = 9500
$ws.Range("P7").Value = 528

# Row 8
$ws.Range("D8").Value = 44804
$ws.Range("J8").Value = 50
$ws.Range("K8").Value = 9500
$ws.Range("L8").Value = 10000
$ws.Range("M8").Value = 9750
$ws.Range("P8").Value = 542
